$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1 "Prix Spot": add a new day column BU (25-aug) with 24 hourly prices
# ---------------------------------------------------------------------------
$wsPrix = $wb.Worksheets.Item("Prix Spot")

# Header cell BU1 — copy the look (bold/border/centered) of the previous
# header cell BT1, then set the label text.
$wsPrix.Range("BT1").Copy()
$wsPrix.Range("BU1").PasteSpecial(-4122)
$wsPrix.Range("BU1").Value = "25-aug"

# Hourly values for column BU, rows 2..25 (row 2 = "00 - 01" ... row 25 = "23 - 24")
$prixValues = @(95.28, 87.03, 74.55, 60.88, 62.43, 80.16, 91.90000000000001, 102.52, 101.82, 91.97, 73.5, 35.02, 10, 5.93, 6.78, 25.2, 56.43, 76.19, 97.43000000000001, 114.78, 125.95, 126.28, 115.04, 105.15)
for ($i = 0; $i -lt $prixValues.Length; $i++) {
    $row = $i + 2
    $wsPrix.Cells.Item($row, 73).Value = $prixValues[$i]
}

# ---------------------------------------------------------------------------
# Sheet 2 "Gaz": append two new daily rows (70, 71)
# ---------------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")

$gazDates = @("2025-08-23", "2025-08-24")
$gazValues = @(32.2, 32.2)
for ($i = 0; $i -lt $gazDates.Length; $i++) {
    $row = 70 + $i
    # Force the date column to Text so the "yyyy-mm-dd" string is kept as
    # literal text instead of being auto-converted to a date serial number.
    $wsGaz.Cells.Item($row, 1).NumberFormat = "@"
    $wsGaz.Cells.Item($row, 1).Value = $gazDates[$i]
    # Re-normalize the cell style back to the plain (unstyled) look used by
    # the rest of the date column (copy format from the last original row).
    $wsGaz.Cells.Item(69, 1).Copy()
    $wsGaz.Cells.Item($row, 1).PasteSpecial(-4122)

    $wsGaz.Cells.Item($row, 2).Value = $gazValues[$i]
}

# ---------------------------------------------------------------------------
# Sheet 3 "CO2": append the same two new daily rows (70, 71); price column
# stays blank, matching the existing blank cell in row 69.
# ---------------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")

$co2Dates = @("2025-08-23", "2025-08-24")
for ($i = 0; $i -lt $co2Dates.Length; $i++) {
    $row = 70 + $i
    $wsCo2.Cells.Item($row, 1).NumberFormat = "@"
    $wsCo2.Cells.Item($row, 1).Value = $co2Dates[$i]
    $wsCo2.Cells.Item(69, 1).Copy()
    $wsCo2.Cells.Item($row, 1).PasteSpecial(-4122)

    # column B is left empty, same as B69 in the source sheet
}

Write-Output "done"
